$d = $word.ActiveDocument

# The page footer ("Ver no Jupiter Salvar em pdf Salvar em docx" and the
# "© 2020 ... Creative Commons Attribution" copyright line), together with
# the blank paragraph that precedes them, is being dropped from the end of
# the document. Locate those paragraphs by content so the edit is resilient
# to any paragraph-numbering differences, then delete the run from the
# blank paragraph right after "LOB1012: ..." through the end of the
# copyright paragraph (leaving the final blank paragraph + page-break
# paragraph untouched).
$count = $d.Paragraphs.Count

$verIdx = -1
$copyIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("Ver no Jupiter")) {
        $verIdx = $i
    }
    if ($t.Contains("Contact:") -and $t.Contains("Creative Commons")) {
        $copyIdx = $i
    }
}

if ($verIdx -gt 0 -and $copyIdx -ge $verIdx) {
    # The paragraph immediately before "Ver no Jupiter ..." is the blank
    # separator paragraph that also needs to go.
    $startPara = $d.Paragraphs.Item($verIdx - 1)
    $endPara = $d.Paragraphs.Item($copyIdx)
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
